# Auto-generated Excel COM-interop edit script
# Applies numeric corrections to H:N columns across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 382.1111
$ws.Range("I33").Value = 382.1111
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 382.1111
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -153.1111
$ws.Range("N33").ClearContents()

$ws.Range("H45").Value = 1000
$ws.Range("J45").Value = 1000
$ws.Range("L45").Value = 3000
$ws.Range("N45").Value = -3384

$ws.Range("H98").Value = 929.3333
$ws.Range("I98").Value = 702.9167
$ws.Range("J98").Value = 1835
$ws.Range("K98").Value = 702.9167
$ws.Range("L98").Value = 1835
$ws.Range("M98").Value = 795.0833
$ws.Range("N98").Value = -4831

$ws.Range("H116").Value = 5338.077
$ws.Range("I116").Value = 2999.8
$ws.Range("J116").Value = 6799.5
$ws.Range("K116").Value = 2999.8
$ws.Range("L116").Value = 6799.5
$ws.Range("M116").Value = 442.1999999999998
$ws.Range("N116").Value = -13683.5

$ws.Range("H122").Value = 929.3333
$ws.Range("I122").Value = 702.9167
$ws.Range("J122").Value = 1835
$ws.Range("K122").Value = 2108.7501
$ws.Range("L122").Value = 5505
$ws.Range("M122").Value = 341.2498999999998
$ws.Range("N122").Value = -10405

$ws.Range("H129").Value = 182793.69
$ws.Range("J129").Value = 197105.02
$ws.Range("L129").Value = 591315.0599999999
$ws.Range("N129").Value = -601315.0599999999

$ws.Range("H132").Value = 2573.8647
$ws.Range("I132").Value = 2603.5278
$ws.Range("K132").Value = 7810.5834
$ws.Range("M132").Value = -5280.5834

$ws.Range("H138").Value = 2426.5312
$ws.Range("J138").Value = 3446.0667
$ws.Range("L138").Value = 10338.2001
$ws.Range("N138").Value = -20618.2001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2841.57
$ws.Range("I32").Value = 2401.7283
$ws.Range("J32").Value = 7899.75
$ws.Range("K32").Value = 2401.7283
$ws.Range("L32").Value = 7899.75
$ws.Range("M32").Value = -2114.7283
$ws.Range("N32").Value = -8473.75

$ws.Range("H45").Value = 2470.5186
$ws.Range("I45").Value = 1906.4615
$ws.Range("K45").Value = 1906.4615
$ws.Range("M45").Value = -1529.4615

$ws.Range("H63").Value = 2950
$ws.Range("I63").Value = 2950
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 2950
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -2264
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 2950
$ws.Range("I66").Value = 2950
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 14750
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -11318
$ws.Range("N66").ClearContents()

$ws.Range("H88").Value = 500977
$ws.Range("J88").Value = 500977
$ws.Range("L88").Value = 500977
$ws.Range("N88").Value = -501789

$ws.Range("H91").Value = 500977
$ws.Range("J91").Value = 500977
$ws.Range("L91").Value = 500977
$ws.Range("N91").Value = -503785

$ws.Range("H122").Value = 1319
$ws.Range("I122").Value = 924.4545000000001
$ws.Range("J122").Value = 9999
$ws.Range("K122").Value = 2773.3635
$ws.Range("L122").Value = 29997
$ws.Range("M122").Value = -323.3635000000004
$ws.Range("N122").Value = -34897

$ws.Range("H132").Value = 14027.214
$ws.Range("I132").Value = 2021.2858
$ws.Range("J132").Value = 74056.86
$ws.Range("K132").Value = 6063.857400000001
$ws.Range("L132").Value = 222170.58
$ws.Range("M132").Value = -3533.857400000001
$ws.Range("N132").Value = -227230.58

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 696.4
$ws.Range("I22").Value = 644.6
$ws.Range("K22").Value = 644.6
$ws.Range("M22").Value = -471.6

$ws.Range("H94").Value = 753.087
$ws.Range("I94").Value = 630.0769
$ws.Range("J94").Value = 913
$ws.Range("K94").Value = 630.0769
$ws.Range("L94").Value = 913
$ws.Range("M94").Value = -179.0769
$ws.Range("N94").Value = -1815

$ws.Range("H99").Value = 770.8333
$ws.Range("I99").Value = 805.3570999999999
$ws.Range("J99").Value = 650
$ws.Range("K99").Value = 805.3570999999999
$ws.Range("L99").Value = 650
$ws.Range("M99").Value = 692.6429000000001
$ws.Range("N99").Value = -3646

$ws.Range("H134").Value = 2941.45
$ws.Range("I134").Value = 2920
$ws.Range("J134").Value = 3349
$ws.Range("K134").Value = 8760
$ws.Range("L134").Value = 10047
$ws.Range("M134").Value = -6225
$ws.Range("N134").Value = -15117

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3747.7144
$ws.Range("J31").Value = 4559.353
$ws.Range("L31").Value = 4559.353
$ws.Range("N31").Value = -5149.353

$ws.Range("H34").Value = 3747.7144
$ws.Range("J34").Value = 4559.353
$ws.Range("L34").Value = 4559.353
$ws.Range("N34").Value = -4963.353

$ws.Range("H132").Value = 4278.1875
$ws.Range("I132").Value = 3037.3635
$ws.Range("J132").Value = 7008
$ws.Range("K132").Value = 9112.0905
$ws.Range("L132").Value = 21024
$ws.Range("M132").Value = -6582.0905
$ws.Range("N132").Value = -26084

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 58.6875
$ws.Range("J12").Value = 97
$ws.Range("L12").Value = 291
$ws.Range("N12").Value = -637

$ws.Range("H37").Value = 62537500
$ws.Range("J37").Value = 62537500
$ws.Range("L37").Value = 187612500
$ws.Range("N37").Value = -187612724

$ws.Range("H68").Value = 1830.4286
$ws.Range("I68").Value = 2000
$ws.Range("J68").Value = 1802.1666
$ws.Range("K68").Value = 6000
$ws.Range("L68").Value = 5406.4998
$ws.Range("M68").Value = -5189
$ws.Range("N68").Value = -7028.4998

$ws.Range("H71").Value = 1830.4286
$ws.Range("I71").Value = 2000
$ws.Range("J71").Value = 1802.1666
$ws.Range("K71").Value = 18000
$ws.Range("L71").Value = 16219.4994
$ws.Range("M71").Value = -13944
$ws.Range("N71").Value = -24331.4994

$ws.Range("H80").Value = 3062
$ws.Range("I80").Value = 2049.5
$ws.Range("J80").Value = 3399.5
$ws.Range("K80").Value = 6148.5
$ws.Range("L80").Value = 10198.5
$ws.Range("M80").Value = -5212.5
$ws.Range("N80").Value = -12070.5

$ws.Range("H83").Value = 3062
$ws.Range("I83").Value = 2049.5
$ws.Range("J83").Value = 3399.5
$ws.Range("K83").Value = 18445.5
$ws.Range("L83").Value = 30595.5
$ws.Range("M83").Value = -13765.5
$ws.Range("N83").Value = -39955.5

$ws.Range("H122").Value = 818.0833
$ws.Range("I122").Value = 470
$ws.Range("J122").Value = 849.7273
$ws.Range("K122").Value = 4230
$ws.Range("L122").Value = 7647.545700000001
$ws.Range("M122").Value = -1780
$ws.Range("N122").Value = -12547.5457

$ws.Range("H129").Value = 197026.34
$ws.Range("I129").Value = 830
$ws.Range("J129").Value = 269309.22
$ws.Range("K129").Value = 2490
$ws.Range("L129").Value = 807927.6599999999
$ws.Range("M129").Value = 2510
$ws.Range("N129").Value = -817927.6599999999

$ws.Range("H131").Value = 147842.77
$ws.Range("J131").Value = 150034
$ws.Range("L131").Value = 450102
$ws.Range("N131").Value = -460182

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 15626381
$ws.Range("I102").Value = 17858244
$ws.Range("J102").Value = 3345.25
$ws.Range("K102").Value = 17858244
$ws.Range("L102").Value = 3345.25
$ws.Range("M102").Value = -17856622
$ws.Range("N102").Value = -6589.25

$ws.Range("H122").Value = 70177880
$ws.Range("I122").Value = 23810840
$ws.Range("K122").Value = 71432520
$ws.Range("M122").Value = -71430070

$ws.Range("H132").Value = 14980.683
$ws.Range("I132").Value = 2708.3784
$ws.Range("K132").Value = 8125.135200000001
$ws.Range("M132").Value = -5595.135200000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4584.909
$ws.Range("I40").Value = 3495
$ws.Range("K40").Value = 3495
$ws.Range("M40").Value = -3359

$ws.Range("H122").Value = 679116.75
$ws.Range("I122").Value = 1155961
$ws.Range("J122").Value = 3587.4167
$ws.Range("K122").Value = 3467883
$ws.Range("L122").Value = 10762.2501
$ws.Range("M122").Value = -3465433
$ws.Range("N122").Value = -15662.2501

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1229.9259
$ws.Range("I122").Value = 1100.2106
$ws.Range("J122").Value = 1538
$ws.Range("K122").Value = 3300.6318
$ws.Range("L122").Value = 4614
$ws.Range("M122").Value = -850.6318000000001
$ws.Range("N122").Value = -9514

$ws.Range("H123").Value = 28607.25
$ws.Range("J123").Value = 28607.25
$ws.Range("L123").Value = 28607.25
$ws.Range("N123").Value = -38407.25

$ws.Range("H126").Value = 1464.5217
$ws.Range("I126").Value = 1477.3846
$ws.Range("J126").Value = 1447.8
$ws.Range("K126").Value = 4432.1538
$ws.Range("L126").Value = 4343.4
$ws.Range("M126").Value = -1962.1538
$ws.Range("N126").Value = -9283.4

$ws.Range("H132").Value = 1540.5385
$ws.Range("I132").Value = 1008.2105
$ws.Range("J132").Value = 2985.4285
$ws.Range("K132").Value = 3024.6315
$ws.Range("L132").Value = 8956.2855
$ws.Range("M132").Value = -494.6315
$ws.Range("N132").Value = -14016.2855

$ws.Range("H136").Value = 28677142
$ws.Range("I136").Value = 41291820
$ws.Range("J136").Value = 7419.091
$ws.Range("K136").Value = 123875460
$ws.Range("L136").Value = 22257.273
$ws.Range("M136").Value = -123872910
$ws.Range("N136").Value = -27357.273

$ws.Range("H140").Value = 30619.75
$ws.Range("J140").Value = 30619.75
$ws.Range("L140").Value = 30619.75
$ws.Range("N140").Value = -40979.75
